# Apply commit: [이용섭] Add - [Table] ItemTable, BasePathTable 수정, 아이템 에셋 - 테이블 바인딩
# Rewrites the "BasePath_BP_File" sheet (2nd sheet) so it becomes an extended
# directory table (Id / Directory) with many more path rows, replacing the
# previous 3-column (Id / Directory_Table_Id / BP_File_Name) mapping table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Clear out the old 3-column data (A1:C4) completely, including the now
# unused column C.
$ws.Cells.Clear()

$data = @(
    @(101, "Level"),
    @(102, "Level/BaseLayerLevel"),
    @(201, "Model"),
    @(301, "TableData"),
    @(401, "UI"),
    @(402, "UI/Widget"),
    @(403, "UI/Widget/LoadingWidget"),
    @(404, "UI/Widget/Practice"),
    @(405, "UI/Widget/Logo"),
    @(406, "UI/Widget/Account"),
    @(407, "UI/Widget/Town/Lobby"),
    @(408, "UI/Widget/Town"),
    @(409, "UI/Widget/Market"),
    @(1001, "UI/Image/Icon"),
    @(5001, "3D/StaticMesh")
)

# Header row (column names) - matches the BasePath_Directory sheet headers.
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "Directory"

# Type row (column data types).
$ws.Range("A2").Value = "int32"
$ws.Range("B2").Value = "int32"

$row = 3
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$ws.Activate()
$ws.Range("A17:B17").Select()

$wb.Save()
